$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Admin Edit User" (row 40) and "Admin Change User Password" (row 41) are now
# implemented, same as the already-checked "Admin Delete User" (row 42) below
# them -> mark the Yes/No column for all three rows as "Yes".
$ws.Range("C40").Value = "Yes"
$ws.Range("C41").Value = "Yes"
$ws.Range("C42").Value = "Yes"

# Reflect where the user ended up scrolled/selected in the sheet afterwards.
$excel.ActiveWindow.ScrollRow = 31
$ws.Range("C43").Select()
